# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows for "Durazno" (Macroferia Regional de Talca)
# right before the existing row 474, shifting the rest of the table down by
# 5 rows (old A1:T537 -> new A1:T542).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 474; everything from the old row 474
# onward shifts down to row 479 onward, unchanged.
$ws.Range("A474:A478").EntireRow.Insert()

# New row 474: Andross / Especial
$ws.Cells.Item(474, 1).Value  = 5
$ws.Cells.Item(474, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(474, 3).Value  = "Maule"
$ws.Cells.Item(474, 4).Value  = 44946
$ws.Cells.Item(474, 5).Value  = 7
$ws.Cells.Item(474, 6).Value  = "Fruta"
$ws.Cells.Item(474, 7).Value  = 100103
$ws.Cells.Item(474, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(474, 9).Value  = 100103004
$ws.Cells.Item(474, 10).Value = "Durazno"
$ws.Cells.Item(474, 11).Value = "Andross"
$ws.Cells.Item(474, 12).Value = "Especial"
$ws.Cells.Item(474, 13).Value = 200
$ws.Cells.Item(474, 14).Value = 16000
$ws.Cells.Item(474, 15).Value = 16000
$ws.Cells.Item(474, 16).Value = 16000
$ws.Cells.Item(474, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(474, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(474, 19).Value = 1067
$ws.Cells.Item(474, 20).Value = 15

# New row 475: Andross / Extra (doble especial)
$ws.Cells.Item(475, 1).Value  = 5
$ws.Cells.Item(475, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(475, 3).Value  = "Maule"
$ws.Cells.Item(475, 4).Value  = 44946
$ws.Cells.Item(475, 5).Value  = 7
$ws.Cells.Item(475, 6).Value  = "Fruta"
$ws.Cells.Item(475, 7).Value  = 100103
$ws.Cells.Item(475, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(475, 9).Value  = 100103004
$ws.Cells.Item(475, 10).Value = "Durazno"
$ws.Cells.Item(475, 11).Value = "Andross"
$ws.Cells.Item(475, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(475, 13).Value = 150
$ws.Cells.Item(475, 14).Value = 17000
$ws.Cells.Item(475, 15).Value = 17000
$ws.Cells.Item(475, 16).Value = 17000
$ws.Cells.Item(475, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(475, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(475, 19).Value = 1133
$ws.Cells.Item(475, 20).Value = 15

# New row 476: Andross / Primera
$ws.Cells.Item(476, 1).Value  = 5
$ws.Cells.Item(476, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(476, 3).Value  = "Maule"
$ws.Cells.Item(476, 4).Value  = 44946
$ws.Cells.Item(476, 5).Value  = 7
$ws.Cells.Item(476, 6).Value  = "Fruta"
$ws.Cells.Item(476, 7).Value  = 100103
$ws.Cells.Item(476, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(476, 9).Value  = 100103004
$ws.Cells.Item(476, 10).Value = "Durazno"
$ws.Cells.Item(476, 11).Value = "Andross"
$ws.Cells.Item(476, 12).Value = "Primera"
$ws.Cells.Item(476, 13).Value = 250
$ws.Cells.Item(476, 14).Value = 15000
$ws.Cells.Item(476, 15).Value = 15000
$ws.Cells.Item(476, 16).Value = 15000
$ws.Cells.Item(476, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(476, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(476, 19).Value = 1000
$ws.Cells.Item(476, 20).Value = 15

# New row 477: Carson / Especial
$ws.Cells.Item(477, 1).Value  = 5
$ws.Cells.Item(477, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(477, 3).Value  = "Maule"
$ws.Cells.Item(477, 4).Value  = 44946
$ws.Cells.Item(477, 5).Value  = 7
$ws.Cells.Item(477, 6).Value  = "Fruta"
$ws.Cells.Item(477, 7).Value  = 100103
$ws.Cells.Item(477, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(477, 9).Value  = 100103004
$ws.Cells.Item(477, 10).Value = "Durazno"
$ws.Cells.Item(477, 11).Value = "Carson"
$ws.Cells.Item(477, 12).Value = "Especial"
$ws.Cells.Item(477, 13).Value = 250
$ws.Cells.Item(477, 14).Value = 16000
$ws.Cells.Item(477, 15).Value = 16000
$ws.Cells.Item(477, 16).Value = 16000
$ws.Cells.Item(477, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(477, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(477, 19).Value = 1067
$ws.Cells.Item(477, 20).Value = 15

# New row 478: Carson / Primera
$ws.Cells.Item(478, 1).Value  = 5
$ws.Cells.Item(478, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(478, 3).Value  = "Maule"
$ws.Cells.Item(478, 4).Value  = 44946
$ws.Cells.Item(478, 5).Value  = 7
$ws.Cells.Item(478, 6).Value  = "Fruta"
$ws.Cells.Item(478, 7).Value  = 100103
$ws.Cells.Item(478, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(478, 9).Value  = 100103004
$ws.Cells.Item(478, 10).Value = "Durazno"
$ws.Cells.Item(478, 11).Value = "Carson"
$ws.Cells.Item(478, 12).Value = "Primera"
$ws.Cells.Item(478, 13).Value = 2003
$ws.Cells.Item(478, 14).Value = 14000
$ws.Cells.Item(478, 15).Value = 14000
$ws.Cells.Item(478, 16).Value = 14000
$ws.Cells.Item(478, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(478, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(478, 19).Value = 933
$ws.Cells.Item(478, 20).Value = 15
